# -----------------------------------------------------------------------
# Ajuste de la Queratina
# % a Marinela
# Menor % a los profesionales de Tocador
#
# The sales-report sheet gets 4 new columns inserted in between the old
# "Valor_producto" / "Part_profesional" columns (Porc_trans, Cost_trans,
# Porc_producto, Valor_Neto), one booking row gets corrected (the
# "Blower cabello largo" line for Sandra Giraldo turns out to really be a
# "Queratina caballero" sale for Carlos Andres Montana, which bumps every
# row below it down by one), the old "descuentos/fondos" column (G) is
# moved to its new home (K), and a couple of new discount/fund rows are
# appended at the bottom.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header row -------------------------------------------------
$headers = @("Fecha de Pago","Nombre cliente","Servicio/Producto","Prestador/Vendedor", `
             "Precio","Porc_trans","Cost_trans","Porc_producto","Valor_producto", `
             "Valor_Neto","Part_profesional","Revisar")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Insert the new Porc_trans / Cost_trans / Porc_producto / Valor_Neto
#        values for the existing booking rows (2-39). Precio (E), and the
#        descriptive columns (A-D) stay as they were; Valor_producto (now
#        column I) and Part_profesional (now column K) keep their original
#        amounts, just recomputed/relocated alongside the new columns.
$rows = @(
    @{Row=2; F=0.036; G=791.9999999999999; H=0.25; I=5500; J=15708; K=8800},
    @{Row=3; F=0.036; G=1008; H=0.26; I=7280; J=19712; K=12600},
    @{Row=4; F=0.036; G=1224; H=0.26; I=8840; J=23936; K=15300},
    @{Row=5; F=0.036; G=4680; H=0.56; I=72800; J=52520; K=10400},
    @{Row=6; F=0.036; G=540; H=0.26; I=3900; J=10560; K=6750},
    @{Row=7; F=0.036; G=1152; H=0.26; I=8320; J=22528; K=14400},
    @{Row=8; F=0; G=0; H=0.26; I=8320; J=23680; K=14400},
    @{Row=9; F=0.036; G=180; H=$null; I=0; J=4820; K=4820},
    @{Row=10; F=0.036; G=1620; H=0.25; I=11250; J=32130; K=18000},
    @{Row=11; F=0.036; G=2160; H=0.25; I=15000; J=42840; K=24000},
    @{Row=12; F=0.036; G=7883.999999999999; H=0.56; I=122640; J=88475.99999999999; K=17520},
    @{Row=13; F=0.036; G=4392; H=0.56; I=68320; J=49288; K=9760},
    @{Row=14; F=0.036; G=3420; H=0.25; I=23750; J=67830; K=38000},
    @{Row=15; F=0.036; G=1620; H=0.25; I=11250; J=32130; K=18000},
    @{Row=16; F=0.036; G=791.9999999999999; H=0.25; I=5500; J=15708; K=8800},
    @{Row=17; F=0.036; G=1620; H=0.25; I=11250; J=32130; K=18000},
    @{Row=18; F=0.036; G=1368; H=0.26; I=9880; J=26752; K=17100},
    @{Row=19; F=0.036; G=1980; H=0.25; I=13750; J=39270; K=22000},
    @{Row=20; F=0.036; G=791.9999999999999; H=0.25; I=5500; J=15708; K=8800},
    @{Row=21; F=0.036; G=1368; H=0.26; I=9880; J=26752; K=17100},
    @{Row=22; F=0.036; G=1368; H=0.26; I=9880; J=26752; K=17100},
    @{Row=23; F=0.036; G=288; H=0.56; I=4480; J=3232; K=640},
    @{Row=24; F=0.036; G=791.9999999999999; H=0.25; I=5500; J=15708; K=8800},
    @{Row=25; F=0; G=0; H=0.25; I=26250; J=78750; K=42000},
    @{Row=26; F=0.036; G=2232; H=0.25; I=15500; J=44268; K=24800},
    @{Row=27; F=0.036; G=1224; H=0.26; I=8840; J=23936; K=15300},
    @{Row=28; F=0.036; G=2160; H=0.25; I=15000; J=42840; K=24000},
    @{Row=29; F=0.036; G=1620; H=0.25; I=11250; J=32130; K=18000},
    @{Row=30; F=0; G=0; H=0.25; I=5500; J=16500; K=8800},
    @{Row=31; F=0.036; G=1980; H=0.25; I=13750; J=39270; K=22000},
    @{Row=32; F=0; G=0; H=0.26; I=15600; J=0; K=-15600},
    @{Row=33; F=0.036; G=720; H=0.25; I=5000; J=14280; K=8000},
    @{Row=34; F=0.036; G=1152; H=0.26; I=8320; J=22528; K=14400},
    @{Row=35; F=0.036; G=306; H=$null; I=0; J=8194; K=8194},
    @{Row=36; F=0.036; G=1260; H=0.1166; I=4081; J=29659; K=19250},
    @{Row=37; F=0.036; G=1800; H=0.25; I=12500; J=35700; K=20000},
    @{Row=38; F=0; G=0; H=0.25; I=3750; J=11250; K=6000},
    @{Row=39; F=0; G=0; H=0.25; I=3750; J=11250; K=6000}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 6).Value = $r.F    # Porc_trans
    $ws.Cells.Item($r.Row, 7).Value = $r.G    # Cost_trans
    if ($null -ne $r.H) {
        $ws.Cells.Item($r.Row, 8).Value = $r.H   # Porc_producto
    }
    $ws.Cells.Item($r.Row, 9).Value = $r.I    # Valor_producto
    $ws.Cells.Item($r.Row, 10).Value = $r.J   # Valor_Neto
    $ws.Cells.Item($r.Row, 11).Value = $r.K   # Part_profesional
}

# --- 3. Row 40 used to be "Blower cabello largo" (Sandra Giraldo); it is
#        really a duplicate "Shampoo Dirigido" for her (A/B/D stay the same).
$ws.Cells.Item(40, 3).Value = "Shampoo Dirigido"
$ws.Cells.Item(40, 5).Value = 22000
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 0
$ws.Cells.Item(40, 8).Value = 0.25
$ws.Cells.Item(40, 9).Value = 5500
$ws.Cells.Item(40, 10).Value = 16500
$ws.Cells.Item(40, 11).Value = 8800

# --- 4. Row 41 used to be Sandra Giraldo's "Shampoo Dirigido"; it now holds
#        the sale that used to sit in row 42 (Katerine Rengifo).
$ws.Cells.Item(41, 1).Value = "01/11/2024 10:50"
$ws.Cells.Item(41, 2).Value = "Katerine Rengifo"
$ws.Cells.Item(41, 3).Value = "Spa Pedicure Tradicional"
$ws.Cells.Item(41, 5).Value = 34000
$ws.Cells.Item(41, 6).Value = 0.036
$ws.Cells.Item(41, 7).Value = 1224
$ws.Cells.Item(41, 8).Value = 0.26
$ws.Cells.Item(41, 9).Value = 8840
$ws.Cells.Item(41, 10).Value = 23936
$ws.Cells.Item(41, 11).Value = 15300

# --- 5. Row 42 becomes the newly-found Queratina sale (Carlos Andres
#        Montana), correcting/replacing the old Katerine Rengifo line.
$ws.Cells.Item(42, 1).Value = "05/11/2024 20:46"
$ws.Cells.Item(42, 2).Value = "Carlos Andres Montana"
$ws.Cells.Item(42, 3).Value = "Queratina caballero - Desde"
$ws.Cells.Item(42, 5).Value = 300000
$ws.Cells.Item(42, 6).Value = 0.036
$ws.Cells.Item(42, 7).Value = 10800
$ws.Cells.Item(42, 8).Value = 0.2016
$ws.Cells.Item(42, 9).Value = 60480
$ws.Cells.Item(42, 10).Value = 228720
$ws.Cells.Item(42, 11).Value = 20000

# --- 6. Rows 43-53 are the Descuentos/Fondos lines for Marinela: their only
#        value column (old G = Part_profesional/"Descuento") moves to the
#        new column K, leaving the old G cell empty.
$discountValues = @{
    43 = -47500
    44 = -20000
    45 = -100000
    46 = -40000
    47 = -17000
    48 = -20000
    49 = -40000
    50 = -31000
    51 = -100000
    52 = -20000
    53 = -115254
}
foreach ($r in $discountValues.Keys) {
    $ws.Cells.Item($r, 11).Value = $discountValues[$r]   # column K
    $ws.Cells.Item($r, 7).Value = ""                     # clear old column G
}

# --- 7. Row 54 used to be a generic "Descuento - Anticipo" line; it is
#        really a "Fondo - Ahorro" entry dated 2024-11-15 for -50000.
#        Force it to stay text so Excel doesn't turn the ISO date into a
#        serial date number.
$ws.Cells.Item(54, 1).NumberFormat = "@"
$ws.Cells.Item(54, 1).Value = "2024-11-15"
$ws.Cells.Item(54, 1).Style = "Normal"
$ws.Cells.Item(54, 3).Value = "Fondo - Ahorro - NA"
$ws.Cells.Item(54, 11).Value = -50000
$ws.Cells.Item(54, 7).Value = ""

# --- 8. New row 55: an additional "Descuento - Producto" line.
$ws.Cells.Item(55, 1).NumberFormat = "@"
$ws.Cells.Item(55, 1).Value = "2024-10-05"
$ws.Cells.Item(55, 1).Style = "Normal"
$ws.Cells.Item(55, 3).Value = "Descuento - Producto - Producto Semi"
$ws.Cells.Item(55, 4).Value = "Marinela Olaya"
$ws.Cells.Item(55, 11).Value = -15600

Write-Host "Sheet updated: new dimension should be A1:L55"
